$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "features_1_2sec" (1st sheet): it was the active tab before, it loses
# that status; its column A gets an explicit width and the whole table gets
# selected.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("features_1_2sec")
$ws1.Select()
$ws1.Columns.Item(1).ColumnWidth = 7.83
$ws1.Range("A1:I7").Select()

# ---------------------------------------------------------------------------
# Sheet "features_2_2sec" (2nd sheet): a few recall numbers were corrected,
# the active cell moves to I3 and column B:I gets reset to the default width.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("features_2_2sec")
$ws2.Select()

$ws2.Range("E4").Value = 0.44
$ws2.Range("E5").Value = 0.25
$ws2.Range("E6").Value = 0.73

$ws2.Range("B1:I1").ColumnWidth = 8.25
$ws2.Range("I3").Select()

# ---------------------------------------------------------------------------
# Sheet "features_1_5sec" (3rd sheet): column A gets an explicit width and
# the active cell moves to H3.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("features_1_5sec")
$ws3.Select()
$ws3.Columns.Item(1).ColumnWidth = 7.83
$ws3.Range("H3").Select()

# ---------------------------------------------------------------------------
# Sheet "features_2_5sec" (4th sheet): becomes the active / selected tab,
# column A gets an explicit width and the active cell moves to A2.
# This sheet is touched last so it remains the active tab on save.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("features_2_5sec")
$ws4.Select()
$ws4.Columns.Item(1).ColumnWidth = 7.83
$ws4.Range("A2").Select()
